$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -0.04699341499341499
$ws.Range("F4").Value = -0.007132596059137766
$ws.Range("G4").Value = 0.02335218640538148
$ws.Range("H4").Value = 0.02695465495465496
$ws.Range("J4").Value = 0.01806486017324933
$ws.Range("C5").Value = -0.02337817137817138
$ws.Range("F5").Value = 0.01359462292114856
$ws.Range("G5").Value = 0.002071828897559966
$ws.Range("H5").Value = 0.006393798393798395
$ws.Range("J5").Value = 0.05561674995045045
$ws.Range("C6").Value = 0.05663590463590464
$ws.Range("F6").Value = -0.01901164495354718
$ws.Range("G6").Value = 0.007710961273337992
$ws.Range("H6").Value = 0.02061980061980062
$ws.Range("J6").Value = 0.02344607958675607
$ws.Range("C7").Value = -0.03400597000597001
$ws.Range("F7").Value = -0.06612641193004069
$ws.Range("G7").Value = 0.03669257042711203
$ws.Range("H7").Value = -0.06378377178377179
$ws.Range("J7").Value = -0.04567677995084063
$ws.Range("C8").Value = 0.2038863598863599
$ws.Range("F8").Value = -0.05614550834892823
$ws.Range("G8").Value = 0.06374238481195806
$ws.Range("H8").Value = 0.9768122208122209
$ws.Range("J8").Value = -0.2388229382558759
$ws.Range("C9").Value = 0.003613263613263614
$ws.Range("F9").Value = 0.03638227887271094
$ws.Range("G9").Value = -0.0444684338717683
$ws.Range("H9").Value = 0.02944744144744145
$ws.Range("J9").Value = -0.02754427570954137
$ws.Range("C10").Value = -0.02215553815553815
$ws.Range("F10").Value = 0.01093921163281705
$ws.Range("G10").Value = -0.02416515715852441
$ws.Range("H10").Value = -0.06988746988746988
$ws.Range("J10").Value = 0.06718400558710962
$ws.Range("C11").Value = -0.04248714648714649
$ws.Range("F11").Value = 0.009315560591769471
$ws.Range("G11").Value = -0.0110987109532978
$ws.Range("H11").Value = -0.01074139074139074
$ws.Range("J11").Value = 0.01008474414525261
$ws.Range("C12").Value = -0.02495224895224895
$ws.Range("F12").Value = -0.01009157050366455
$ws.Range("G12").Value = 0.0160169468422559
$ws.Range("H12").Value = -0.05575305175305176
$ws.Range("J12").Value = -0.02074856687305827
$ws.Range("C13").Value = 0.01924258324258324
$ws.Range("F13").Value = -0.03423690540257183
$ws.Range("G13").Value = 0.004653890044915938
$ws.Range("H13").Value = 0.01704110904110904
$ws.Range("J13").Value = 0.0005873765909008505
$ws.Range("C14").Value = 0.8855650775650776
$ws.Range("F14").Value = 0.01629107299092888
$ws.Range("G14").Value = -0.009668677279661132
$ws.Range("H14").Value = 0.0002385002385002385
$ws.Range("J14").Value = 0.01432310040903901
$ws.Range("C15").Value = -0.01503491103491103
$ws.Range("F15").Value = -0.03098692648234342
$ws.Range("G15").Value = 0.009519416528067226
$ws.Range("H15").Value = -0.0570926970926971
$ws.Range("J15").Value = 0.001888753900086423
$ws.Range("C16").Value = -0.05344660144660145
$ws.Range("F16").Value = 0.02241836154787277
$ws.Range("G16").Value = -0.03129585673302161
$ws.Range("H16").Value = -0.002793350793350793
$ws.Range("J16").Value = -0.04088310512840376
$ws.Range("C17").Value = -0.0005709605709605709
$ws.Range("F17").Value = -0.03194182930458597
$ws.Range("G17").Value = 0.03691778672048535
$ws.Range("H17").Value = 0.009723249723249724
$ws.Range("J17").Value = 0.07000457242459987
$ws.Range("C18").Value = 0.1152071712071712
$ws.Range("F18").Value = 0.004248536662775623
$ws.Range("G18").Value = -0.002940054118279753
$ws.Range("H18").Value = -0.01976607176607177
$ws.Range("J18").Value = 0.01915251126742634
$ws.Range("C19").Value = -0.2180841020841021
$ws.Range("F19").Value = 0.007870492894754628
$ws.Range("G19").Value = 0.0221312193387525
$ws.Range("H19").Value = -0.005386757386757387
$ws.Range("J19").Value = 0.0008608568660220073
$ws.Range("C20").Value = 0.03598840798840799
$ws.Range("F20").Value = -0.05651264185834246
$ws.Range("G20").Value = 0.03522421221017902
$ws.Range("H20").Value = 0.01699283299283299
$ws.Range("J20").Value = -0.03073195891635067
$ws.Range("C21").Value = -0.01506323106323107
$ws.Range("F21").Value = 0.000773806528125535
$ws.Range("G21").Value = -0.002063506359800996
$ws.Range("H21").Value = -0.001079233079233079
$ws.Range("J21").Value = -0.03102705121321352
$ws.Range("C22").Value = -0.09701215301215302
$ws.Range("F22").Value = 0.03408426492969901
$ws.Range("G22").Value = -0.02833595608682705
$ws.Range("H22").Value = 0.02665963465963466
$ws.Range("J22").Value = -0.009219447274763959
$ws.Range("C23").Value = 0.01524050724050724
$ws.Range("F23").Value = -0.007903980686095118
$ws.Range("G23").Value = 0.02718125351168523
$ws.Range("H23").Value = 0.01521103521103521
$ws.Range("J23").Value = 0.01005485811518727
$ws.Range("C24").Value = -0.1142429582429582
$ws.Range("F24").Value = 0.00231388680405151
$ws.Range("G24").Value = 0.01513157496450145
$ws.Range("H24").Value = 0.1370689010689011
$ws.Range("J24").Value = 0.008122472171207006
$ws.Range("C25").Value = 0.02023979623979624
$ws.Range("F25").Value = -0.09443721046067664
$ws.Range("G25").Value = 0.08309410565669759
$ws.Range("H25").Value = -0.03316828516828517
$ws.Range("J25").Value = -0.01853676064798121
$ws.Range("C26").Value = 0.06612943812943814
$ws.Range("F26").Value = 0.00493151279220275
$ws.Range("G26").Value = -0.002904435638227972
$ws.Range("H26").Value = -0.001706197706197706
$ws.Range("J26").Value = 0.004132054156846482
$ws.Range("C27").Value = 0.0465970785970786
$ws.Range("F27").Value = -0.1592702179082092
$ws.Range("G27").Value = 0.2022380328924657
$ws.Range("H27").Value = -0.02251911451911452
$ws.Range("J27").Value = 0.02163654376636303
$ws.Range("C28").Value = 0.004948588948588949
$ws.Range("F28").Value = -0.009002309830901469
$ws.Range("G28").Value = 0.01860453157866024
$ws.Range("H28").Value = -0.03140517140517141
$ws.Range("J28").Value = 0.0430870753455978
$ws.Range("C29").Value = -0.01427672627672628
$ws.Range("F29").Value = 0.04979917430963807
$ws.Range("G29").Value = -0.04973264954470642
$ws.Range("H29").Value = 0.008113916113916115
$ws.Range("J29").Value = 0.04074664099112084
$ws.Range("C30").Value = 0.2272703032703032
$ws.Range("F30").Value = -0.7683386218888937
$ws.Range("G30").Value = 0.9478196497099235
$ws.Range("H30").Value = -0.03066137466137466
$ws.Range("J30").Value = -0.03440390861033206
